$d = $word.ActiveDocument

# Locate the target paragraph: the last (empty) paragraph of the
# bulleted list at the end of the document's body, right before the
# sectPr. We find it robustly as the last paragraph in the document
# whose text is empty.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "") {
        $target = $para
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs.Item($count)
}

# Select the whole paragraph, including its end-of-paragraph mark, so
# that InsertXML replaces the paragraph's content in place (preserving
# its pPr: list style, numbering, justification, run-mark language)
# rather than inserting a sibling paragraph before it.
$r = $target.Range
$null = $r.MoveEnd(1, 1)

$run1 = "En el Branch alterno se elimina "
$enDash = [char]0x2013
$run2 = "CASO 5 DERRATEO DE TRANSFORMADORES " + $enDash + "FACTOR K"
$run3 = ".pdf y se inserta el archivo flujo de potencia.pdf"

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">' + $run1 + '</w:t></w:r>' +
       '<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>' + $run2 + '</w:t></w:r>' +
       '<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>' + $run3 + '</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$null = $r.InsertXML($xml)
